$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 93
$ws1.Range("F20").Value = 296
$ws1.Range("F21").Value = 3246
$ws1.Range("F23").Value = 414
$ws1.Range("F27").Value = 2863
$ws1.Range("F28").Value = 1668
$ws1.Range("F34").Value = 1941
$ws1.Range("F36").Value = 1952
$ws1.Range("F39").Value = 113
$ws1.Range("F42").Value = 917
$ws1.Range("F44").Value = 1061
$ws1.Range("F49").Value = 3386

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 1
$ws2.Range("F8").Value = 16

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 93
$ws4.Range("F20").Value = 296
$ws4.Range("F21").Value = 3246
$ws4.Range("F23").Value = 414
$ws4.Range("F24").Value = 16
$ws4.Range("F26").Value = 2863
$ws4.Range("F27").Value = 1668
$ws4.Range("F33").Value = 1941
$ws4.Range("F36").Value = 1952
$ws4.Range("F38").Value = 113
$ws4.Range("F40").Value = 917
$ws4.Range("F42").Value = 1061
$ws4.Range("F48").Value = 3386
